$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 342
$ws.Range("F2").Value = 45444
$ws.Range("G2").Value = 30319
$ws.Range("H2").Value = 45505
# Row 3
$ws.Range("E3").Value = 30256
$ws.Range("F3").Value = 45444
# Row 4
$ws.Range("E4").Value = 30256
$ws.Range("F4").Value = 45444
$ws.Range("G4").Value = 30319
$ws.Range("H4").Value = 45505
# Row 5
$ws.Range("E5").Value = 30256
$ws.Range("F5").Value = 45444
$ws.Range("G5").Value = 30319
$ws.Range("H5").Value = 45505
# Row 6
$ws.Range("C6").Value = 451
$ws.Range("F6").Value = 45444
$ws.Range("G6").Value = 30319
$ws.Range("H6").Value = 45505
# Row 7
$ws.Range("E7").Value = 30256
$ws.Range("F7").Value = 45444
$ws.Range("G7").Value = 30319
$ws.Range("H7").Value = 45505
# Row 8
$ws.Range("D8").Value = 418
$ws.Range("E8").Value = 30256
$ws.Range("F8").Value = 45444
$ws.Range("H8").Value = 45505
# Row 9
$ws.Range("E9").Value = 30256
$ws.Range("F9").Value = 45444
$ws.Range("G9").Value = 30319
$ws.Range("H9").Value = 45505
# Row 10
$ws.Range("D10").Value = 491
$ws.Range("E10").Value = 30256
$ws.Range("F10").Value = 45444
$ws.Range("H10").Value = 45505
# Row 11
$ws.Range("E11").Value = 30225
$ws.Range("F11").Value = 45413
$ws.Range("G11").Value = 30319
$ws.Range("H11").Value = 45505
# Row 12
$ws.Range("C12").Value = 379
$ws.Range("D12").Value = 360
$ws.Range("F12").Value = 45444
$ws.Range("H12").Value = 45505
# Row 13
$ws.Range("C13").Value = 474
$ws.Range("F13").Value = 45413
$ws.Range("G13").Value = 30319
$ws.Range("H13").Value = 45505
# Row 14
$ws.Range("D14").Value = 404
$ws.Range("H14").Value = 45505
# Row 15
$ws.Range("G15").Value = 30286
$ws.Range("H15").Value = 45505
# Row 16
$ws.Range("C16").Value = 462
$ws.Range("D16").Value = 418
$ws.Range("F16").Value = 45413
$ws.Range("H16").Value = 45505
# Row 17
$ws.Range("C17").Value = 378
$ws.Range("D17").Value = 402
$ws.Range("F17").Value = 45444
$ws.Range("H17").Value = 45505
# Row 18
$ws.Range("D18").Value = 268
$ws.Range("E18").Value = 30256
$ws.Range("F18").Value = 45444
$ws.Range("H18").Value = 45505
# Row 19
$ws.Range("D19").Value = 406
$ws.Range("E19").Value = 30256
$ws.Range("F19").Value = 45444
$ws.Range("H19").Value = 45505
# Row 20
$ws.Range("C20").Value = 489
$ws.Range("F20").Value = 45444
$ws.Range("G20").Value = 30319
$ws.Range("H20").Value = 45505
# Row 21
$ws.Range("E21").Value = 30225
$ws.Range("F21").Value = 45413
$ws.Range("G21").Value = 30319
$ws.Range("H21").Value = 45505
# Row 22
$ws.Range("D22").Value = 388
$ws.Range("E22").Value = 30256
$ws.Range("F22").Value = 45444
$ws.Range("H22").Value = 45505
# Row 23
$ws.Range("D23").Value = 225
$ws.Range("E23").Value = 30225
$ws.Range("F23").Value = 45413
$ws.Range("H23").Value = 45505
# Row 24
$ws.Range("C24").Value = 402
$ws.Range("F24").Value = 45444
$ws.Range("G24").Value = 30319
$ws.Range("H24").Value = 45505
# Row 25
$ws.Range("G25").Value = 30319
$ws.Range("H25").Value = 45505
# Row 26
$ws.Range("C26").Value = 462
$ws.Range("D26").Value = 406
$ws.Range("F26").Value = 45444
$ws.Range("H26").Value = 45505
# Row 27
$ws.Range("D27").Value = 418
$ws.Range("E27").Value = 30256
$ws.Range("F27").Value = 45444
$ws.Range("H27").Value = 45505
# Row 28
$ws.Range("C28").Value = 367
$ws.Range("D28").Value = 375
$ws.Range("F28").Value = 45444
$ws.Range("H28").Value = 45505
